$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65 (pushing ZA5480... etc. down by one),
# and populate it with the newly catalogued Eurobarometer wave so it
# lands as a new shared-string entry.
$ws.Rows.Item(65).Insert()
$ws.Range("A65").Value = "ZA5481: Eurobarometer 75.3 (May 2011) Europe 2020, Financial and Economic Crisis, European Union Budget, and the Common Agricultural Policy"

# Match the author's final selection/scroll position in the sheet.
$ws.Range("A65").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49

# Page setup was touched as part of this edit (portrait orientation).
$ws.PageSetup.Orientation = 1
